$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 29
$ws1.Range("F8").Value = 7941
$ws1.Range("F9").Value = 128
$ws1.Range("F11").Value = 6823
$ws1.Range("F12").Value = 159
$ws1.Range("F13").Value = 295
$ws1.Range("F14").Value = 4855
$ws1.Range("F16").Value = 5266
$ws1.Range("F17").Value = 1067
$ws1.Range("F18").Value = 310
$ws1.Range("F19").Value = 313
$ws1.Range("F20").Value = 427
$ws1.Range("F21").Value = 306
$ws1.Range("F22").Value = 252
$ws1.Range("F23").Value = 135
$ws1.Range("F26").Value = 8951
$ws1.Range("F28").Value = 1594
$ws1.Range("F30").Value = 37
$ws1.Range("F37").Value = 1839
$ws1.Range("F39").Value = 1136
$ws1.Range("F41").Value = 4663
$ws1.Range("F43").Value = 1151
$ws1.Range("F44").Value = 63
$ws1.Range("F45").Value = 138
$ws1.Range("F48").Value = 1225

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F17").Value = 884

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 29
$ws4.Range("F9").Value = 7941
$ws4.Range("F10").Value = 128
$ws4.Range("F12").Value = 6823
$ws4.Range("F13").Value = 159
$ws4.Range("F14").Value = 295
$ws4.Range("F16").Value = 4855
$ws4.Range("F18").Value = 5266
$ws4.Range("F19").Value = 1067
$ws4.Range("F20").Value = 310
$ws4.Range("F21").Value = 313
$ws4.Range("F22").Value = 427
$ws4.Range("F23").Value = 306
$ws4.Range("F24").Value = 252
$ws4.Range("F25").Value = 135
$ws4.Range("F27").Value = 8951
$ws4.Range("F29").Value = 1594
$ws4.Range("F30").Value = 37
$ws4.Range("F32").Value = 817
$ws4.Range("F37").Value = 1839
$ws4.Range("F39").Value = 1136
$ws4.Range("F41").Value = 4663
$ws4.Range("F43").Value = 1151
$ws4.Range("F44").Value = 63
$ws4.Range("F45").Value = 138
$ws4.Range("F48").Value = 1225
